# Weekly update: a new price-sampling row for Jengibre (La Palmera de La
# Serena) is inserted before the current row 52, pushing the existing
# rows 52-160 down by one (they keep their original data) and adding a
# fresh data point at the end of the range (new row 161, which is a copy
# of what used to be row 160 before the shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 52; everything below
# (rows 52-160) shifts down to 53-161.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly observation.
$ws.Cells.Item(52, 1).Value = 8
$ws.Cells.Item(52, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(52, 3).Value = "Coquimbo"
$ws.Cells.Item(52, 4).Value = 45133
$ws.Cells.Item(52, 5).Value = 4
$ws.Cells.Item(52, 6).Value = 100114007
$ws.Cells.Item(52, 7).Value = "Jengibre"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 400
$ws.Cells.Item(52, 11).Value = 17000
$ws.Cells.Item(52, 12).Value = 18000
$ws.Cells.Item(52, 13).Value = 17500
$ws.Cells.Item(52, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(52, 15).Value = "Perú"
$ws.Cells.Item(52, 16).Value = 1346
$ws.Cells.Item(52, 17).Value = 13
$ws.Cells.Item(52, 18).Value = "Hortaliza"
